# Update the "ランサーズ" (Lancers) listing sheet to the newer scrape snapshot
# (2025-10-20 06:26:54 JST): refresh row 2-3 timestamps, replace the
# remaining data rows 4-12 with the new listings, drop rows 13-14, and
# resize columns D and H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = "2025-10-20 06:26:54"

# ---------------------------------------------------------------------
# 1) Drop the two trailing rows (13 & 14) completely - shrinks the used
#    range from A1:H14 down to A1:H12.
# ---------------------------------------------------------------------
$ws.Range("A13:H14").Delete()

# ---------------------------------------------------------------------
# 2) Wipe every existing hyperlink on the sheet. (Hyperlinks.Delete on
#    any scoped range clears the whole collection in this engine, so one
#    call is enough - we rebuild the 11 that should remain below.)
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 6).Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 3) Refresh the "取得日時" timestamp for every remaining row (2-12).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
  $ws.Cells.Item($r, 1).Value = $newDate
}

# ---------------------------------------------------------------------
# 4) Rows 2 & 3 keep their original listing data (only the timestamp
#    changed, handled above). Rows 4-12 get entirely new listing data.
# ---------------------------------------------------------------------

# Row 4
$ws.Cells.Item(4, 2).Value = "Kintone × SharePoint × PowerAutomate連携業務システム(AI連携)"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "5,000,000 円 ~ / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5416528"
$ws.Cells.Item(4, 7).Value = 325
$ws.Cells.Item(4, 8).Value = "🔥AI,Ai"

# Row 5
$ws.Cells.Item(5, 2).Value = "【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5416511"
$ws.Cells.Item(5, 7).Value = 155
$ws.Cells.Item(5, 8).Value = "◆開発,Node.js"

# Row 6
$ws.Cells.Item(6, 2).Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5416508"
$ws.Cells.Item(6, 7).Value = 68
$ws.Cells.Item(6, 8).Value = "◆開発"

# Row 7
$ws.Cells.Item(7, 2).Value = "【急募】全国店舗をGoogleマップで表示するWPプラグイン開発"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5416539"
$ws.Cells.Item(7, 7).Value = 63
$ws.Cells.Item(7, 8).Value = "◆開発"

# Row 8
$ws.Cells.Item(8, 2).Value = "【急募】GASを使った顧客管理スプレッドシートの作成・改修依頼"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5416338"
$ws.Cells.Item(8, 7).Value = 33
$ws.Cells.Item(8, 8).Value = "◇管理"

# Row 9
$ws.Cells.Item(9, 2).Value = "サイトスピードが遅く サイトスピードを速くしたい ワードプレス"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "20,000 円 ~ 30,000 円 / 募集期間 3 日、取引期間 0 日"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5416402"
$ws.Cells.Item(9, 7).Value = 30
$ws.Cells.Item(9, 8).Value = "◇サイト"

# Row 10 (no H value: the skill-summary cell is removed entirely)
$ws.Cells.Item(10, 2).Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Cells.Item(10, 7).Value = 25
$ws.Cells.Item(10, 8).ClearContents()

# Row 11 (no H value)
$ws.Cells.Item(11, 2).Value = "【継続案件あり】AWSに精通しているインフラエンジニアを募集します"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5416510"
$ws.Cells.Item(11, 7).Value = 10
$ws.Cells.Item(11, 8).ClearContents()

# Row 12 (no H value, as before)
$ws.Cells.Item(12, 2).Value = "【急募】エクセルマクロの組み方を教えてください!"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5416433"
$ws.Cells.Item(12, 7).Value = 10
$ws.Cells.Item(12, 8).ClearContents()

# ---------------------------------------------------------------------
# 5) Rebuild the hyperlinks for F2:F12 against the (possibly new) URLs
#    now sitting in those cells.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
  $cell = $ws.Cells.Item($r, 6)
  $ws.Hyperlinks.Add($cell, $cell.Value())
}

# ---------------------------------------------------------------------
# 6) Resize column D (32 -> 41) and column H (19 -> 17). ColumnWidth in
#    this object model reports 5/6 of a character wider than the raw
#    OOXML <col width> value, so subtract that fixed offset to land on
#    the exact target width.
# ---------------------------------------------------------------------
$widthOffset = 5 / 6
$ws.Columns.Item(4).ColumnWidth = 41 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 17 - $widthOffset
